# Apply vendor / component parameter updates to the bici-pcb BOM sheet.
# Values are prefixed with a leading apostrophe so Excel keeps treating
# these as literal text (preserves the existing "quote prefix" cell
# formatting instead of resetting it).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 4 - KGM15AR70J104JT 0.1uF capacitor: one more designator (C406) added,
# quantity bumped 8 -> 9
$ws.Range("C4").Value = "'C6, C101, C201, C202, C402, C403, C404, C405, C406"
$ws.Range("D4").Value = 9

# Row 10 - J1 USB-C receptacle: switched from Amphenol 10164359-00011LF to
# GCT USB4125-GF-A
$ws.Range("A10").Value = "'USB4125-GF-A"
$ws.Range("B10").Value = "'CONN RCPT TYPE C 6P SMD RA"
$ws.Range("E10").Value = "'GCT"
$ws.Range("F10").Value = "'USB4125-GF-A"
$ws.Range("H10").Value = "'2073-USB4125-GF-ACT-ND"

# Row 16 - R3/R401 470 ohm resistor: vendor switched DigiKey -> Mouser
$ws.Range("G16").Value = "'Mouser"
$ws.Range("H16").Value = "'667-ERJ-3EKF4700V"

# Row 21 - R215/R216/R301/R303/R402/R404 10k resistor: vendor switched
# DigiKey -> Mouser
$ws.Range("G21").Value = "'Mouser"
$ws.Range("H21").Value = "'667-ERJ-3EKF1002V"

# Row 26 - U101 IMU (ICM-42670-P): vendor switched DigiKey -> Mouser,
# package label reformatted
$ws.Range("G26").Value = "'Mouser"
$ws.Range("H26").Value = "'410-ICM-42670-P"
$ws.Range("J26").Value = "'LGA-14"
